$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 6 (pushes rows 6-38 down to 7-39), then populate
# the new row's A/B/C/D cells with the "start" data point.
$ws.Rows.Item(6).Insert()

$ws.Range("A6").Value = "start"
$ws.Range("B6").Value = 1
$ws.Range("C6").Value = "start"
$ws.Range("D6").Value = "start"

$ws.Range("D6").Select()
